$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2/K2: a new REPORT_TYPE_CODE "002" row is inserted ahead of the
# existing "001" value, so J2 becomes "002" and K2 keeps "001".
# Force text so the numeric-looking codes aren't coerced to numbers,
# then drop the number-format style again so no stray formatting sticks.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("J2").ClearFormats()

$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "001"
$ws.Range("K2").ClearFormats()

# Report date moves to the newer period
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Updated financial figures for the new reporting period
$ws.Range("O2").Value = 44071294.14
$ws.Range("P2").Value = 488068692.66
$ws.Range("Q2").Value = 428381181.29
$ws.Range("R2").ClearContents()
$ws.Range("S2").Value = 387860071.43
$ws.Range("T2").Value = 387860071.43
$ws.Range("U2").ClearContents()
$ws.Range("V2").Value = 4710234.38
$ws.Range("W2").Value = 32841329.34
$ws.Range("X2").Value = 54502.32
$ws.Range("Y2").Value = 62546378.87
$ws.Range("Z2").Value = 62482519.65
$ws.Range("AA2").Value = 15260027.13
$ws.Range("AG2").Value = 2915043.82
$ws.Range("AP2").Value = 21.8469833909
$ws.Range("AQ2").Value = 57.055420470867
$ws.Range("AR2").Value = 61.308047011112
$ws.Range("AS2").Value = 42115494.14
$ws.Range("AT2").Value = 34.49544448596
